$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 previously used the "highlighted/flagged" style (blue fill, style index 6)
# because the "NO OF HOURS LATE" (column F) value was missing/uncomputed.
# Fix: set F12 to the correctly computed value (1.25) and restore the row's
# formatting to match the normal (unhighlighted) rows, like row 11.
$ws.Range("A11:J11").Copy() | Out-Null
$ws.Range("A12:J12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F12").Value = 1.25
